$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75 (everything from old row 75 down shifts to 76+).
$ws.Rows("75:75").Insert()

# Populate the new "Officer ID" row (now row 75).
$ws.Range("A75").Value = "Officer ID"
$ws.Range("B75").Value = "A unique identifier assigned to an officer"
$ws.Range("C75").Value = "Officer ID"
$ws.Range("D75").Value = "ID34567"

# The former "Officer Badge Number" row (old row 75) is now row 76: update its
# Vermont eCitation Element label (new shared string) before the new row's
# mapping path, matching the shared-string table order of the authored edit.
$ws.Range("C76").Value = "Officer Badge No."

$ws.Range("E75").Value = "/ir-doc:IncidentReport/lexspd:doPublish/lexs:PublishMessageContainer/lexs:PublishMessage/lexs:DataItemPackage/lexs:StructuredPayload/inc-ext:IncidentReport/inc-ext:EnforcementOfficial/inc-ext:EnforcementOfficialIdentification/nc:IdentificationID"

# Match the row height used for this new row in the edited workbook.
$ws.Rows("75:75").RowHeight = 42

$ws.Range("E76").Value = "/ir-doc:IncidentReport/lexspd:doPublish/lexs:PublishMessageContainer/lexs:PublishMessage/lexs:DataItemPackage/lexs:Digest/lexsdigest:EntityPerson/j:EnforcementOfficial[@s:id=/ir-doc:IncidentReport/lexspd:doPublish/lexs:PublishMessageContainer/lexs:PublishMessage/lexs:DataItemPackage/lexs:Digest/lexsdigest:Associations/j:ActivityEnforcementOfficialAssociation/nc:PersonReference/@s:ref]/j:EnforcementOfficialBadgeIdentification/nc:IdentificationID"

# Restore the frozen-pane view state (top-left cell of the scrolled region and
# the active selection) to reflect where the editor was working.
$ws.Range("D76").Select()
